$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the mailto hyperlinks from column A (and their special formatting)
# ---------------------------------------------------------------------------
$ws.Range("A1:A4").ClearFormats()
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Replace the sender addresses in column A with the new client list, and
#    add a 5th row for the newly integrated sender
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "uwtsjgxgxgk42@gmail.com"
$ws.Range("A2").Value = "kdhuhfhfi38@gmail.com"
$ws.Range("A3").Value = "kfhhfhfbc@gmail.com"
$ws.Range("A4").Value = "jyretikdgdhl@gmail.com"
$ws.Range("A5").Value = "guuofkhc@gmail.com"

# ---------------------------------------------------------------------------
# 3. Add the authentication token/password column (B) next to every sender,
#    styled with a small Arial font and a medium grey box border.
#    Build the style once on a scratch cell, then copy/paste the format so
#    every destination cell shares a single stylesheet entry.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z100")
$scratch.Font.Size = 10
$scratch.Font.Name = "Arial"
$scratch.Borders.Weight = -4138
$scratch.Borders.Color = 13421772

$ws.Range("B1:B5").Value = "aass1122"
$scratch.Copy()
$ws.Range("B1:B5").PasteSpecial(-4122)
$scratch.Clear()

# ---------------------------------------------------------------------------
# 4. Give the new rows a little extra height to match the boxed look, and
#    move the active selection the way the author left it.
# ---------------------------------------------------------------------------
$ws.Rows("1:5").RowHeight = 15.75
$ws.Range("B8").Select()
